# Apply odds/score updates to Sheet1 per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Cells.Item(12, 7).Value = 2.7  # G12: 2.45 -> 2.7
$ws.Cells.Item(12, 9).Value = 2.8  # I12: 3.2 -> 2.8
$ws.Cells.Item(12, 10).Value = 3.5  # J12: 3.25 -> 3.5
$ws.Cells.Item(12, 12).Value = 3.75  # L12: 4 -> 3.75
$ws.Cells.Item(12, 17).Value = 1.98  # Q12: 2 -> 1.98
$ws.Cells.Item(12, 18).Value = 1.88  # R12: 1.85 -> 1.88
$ws.Cells.Item(12, 21).Value = 4.2  # U12: 4.3 -> 4.2
$ws.Cells.Item(12, 22).Value = 1.22  # V12: 1.21 -> 1.22
$ws.Cells.Item(12, 23).Value = 5  # W12: 5.5 -> 5
$ws.Cells.Item(12, 24).Value = 1.17  # X12: 1.14 -> 1.17
$ws.Cells.Item(12, 29).Value = 6.5  # AC12: 6 -> 6.5
$ws.Cells.Item(12, 30).Value = 12  # AD12: 10 -> 12
$ws.Cells.Item(12, 32).Value = 26  # AF12: 23 -> 26
$ws.Cells.Item(12, 33).Value = 26  # AG12: 23 -> 26
$ws.Cells.Item(12, 35).Value = 6.5  # AI12: 6 -> 6.5
$ws.Cells.Item(12, 37).Value = 17  # AK12: 19 -> 17
$ws.Cells.Item(12, 40).Value = 7  # AN12: 7.5 -> 7
$ws.Cells.Item(12, 41).Value = 13  # AO12: 15 -> 13
$ws.Cells.Item(12, 42).Value = 12  # AP12: 13 -> 12
$ws.Cells.Item(12, 43).Value = 29  # AQ12: 34 -> 29
$ws.Cells.Item(12, 44).Value = 29  # AR12: 34 -> 29

# Row 20
$ws.Cells.Item(20, 7).Value = 2.47  # G20: 2.5 -> 2.47
$ws.Cells.Item(20, 8).Value = 3.2  # H20: 3.15 -> 3.2
$ws.Cells.Item(20, 9).Value = 2.7  # I20: 2.72 -> 2.7
$ws.Cells.Item(20, 10).Value = 3.15  # J20: 3.1 -> 3.15
$ws.Cells.Item(20, 11).Value = 2.02  # K20: 2 -> 2.02
$ws.Cells.Item(20, 12).Value = 3.3  # L20: 3.4 -> 3.3
$ws.Cells.Item(20, 23).Value = 2.95  # W20: 2.92 -> 2.95
$ws.Cells.Item(20, 26).Value = 2.47  # Z20: 2.45 -> 2.47
$ws.Cells.Item(20, 29).Value = 8.25  # AC20: 8.75 -> 8.25
$ws.Cells.Item(20, 30).Value = 12.5  # AD20: 13 -> 12.5
$ws.Cells.Item(20, 31).Value = 9.5  # AE20: 9.25 -> 9.5
$ws.Cells.Item(20, 32).Value = 27  # AF20: 28 -> 27
$ws.Cells.Item(20, 34).Value = 29  # AH20: 27 -> 29
$ws.Cells.Item(20, 35).Value = 9.75  # AI20: 9.5 -> 9.75
$ws.Cells.Item(20, 36).Value = 6.2  # AJ20: 6.1 -> 6.2
$ws.Cells.Item(20, 37).Value = 13  # AK20: 12.5 -> 13
$ws.Cells.Item(20, 40).Value = 9  # AN20: 8.75 -> 9
$ws.Cells.Item(20, 41).Value = 14.5  # AO20: 14 -> 14.5
$ws.Cells.Item(20, 42).Value = 9.75  # AP20: 10 -> 9.75
$ws.Cells.Item(20, 44).Value = 22  # AR20: 23 -> 22
$ws.Cells.Item(20, 45).Value = 29  # AS20: 30 -> 29

# Row 35
$ws.Cells.Item(35, 7).Value = 1.91  # G35: 1.95 -> 1.91
$ws.Cells.Item(35, 9).Value = 3.75  # I35: 3.6 -> 3.75
$ws.Cells.Item(35, 12).Value = 4.5  # L35: 4.33 -> 4.5
$ws.Cells.Item(35, 13).Value = 1.07  # M35: 1.03 -> 1.07
$ws.Cells.Item(35, 15).Value = 1.36  # O35: 1.3 -> 1.36
$ws.Cells.Item(35, 23).Value = 3.75  # W35: 4 -> 3.75
$ws.Cells.Item(35, 24).Value = 1.25  # X35: 1.16 -> 1.25
$ws.Cells.Item(35, 25).Value = 1.5  # Y35: 1.44 -> 1.5
$ws.Cells.Item(35, 26).Value = 2.5  # Z35: 2.63 -> 2.5
$ws.Cells.Item(35, 27).Value = 1.91  # AA35: 1.87 -> 1.91
$ws.Cells.Item(35, 28).Value = 1.8  # AB35: 1.77 -> 1.8
$ws.Cells.Item(35, 30).Value = 8.5  # AD35: 9 -> 8.5
$ws.Cells.Item(35, 34).Value = 34  # AH35: 29 -> 34
$ws.Cells.Item(35, 35).Value = 8  # AI35: 8.5 -> 8
$ws.Cells.Item(35, 42).Value = 15  # AP35: 13 -> 15

# Row 43
$ws.Cells.Item(43, 7).Value = 1.62  # G43: 1.55 -> 1.62
$ws.Cells.Item(43, 8).Value = 4  # H43: 4.2 -> 4
$ws.Cells.Item(43, 9).Value = 5.25  # I43: 5.75 -> 5.25
$ws.Cells.Item(43, 10).Value = 2.2  # J43: 2.1 -> 2.2
$ws.Cells.Item(43, 11).Value = 2.3  # K43: 2.38 -> 2.3
$ws.Cells.Item(43, 14).Value = 12  # N43: 13 -> 12
$ws.Cells.Item(43, 15).Value = 1.25  # O43: 1.22 -> 1.25
$ws.Cells.Item(43, 16).Value = 3.75  # P43: 4 -> 3.75
$ws.Cells.Item(43, 19).Value = 1.83  # S43: 1.69 -> 1.83
$ws.Cells.Item(43, 20).Value = 2.03  # T43: 2.07 -> 2.03
$ws.Cells.Item(43, 23).Value = 3  # W43: 2.75 -> 3
$ws.Cells.Item(43, 24).Value = 1.36  # X43: 1.4 -> 1.36
$ws.Cells.Item(43, 25).Value = 1.36  # Y43: 1.33 -> 1.36
$ws.Cells.Item(43, 26).Value = 3  # Z43: 3.25 -> 3
$ws.Cells.Item(43, 30).Value = 8  # AD43: 7.5 -> 8
$ws.Cells.Item(43, 32).Value = 12  # AF43: 11 -> 12
$ws.Cells.Item(43, 33).Value = 13  # AG43: 12 -> 13
$ws.Cells.Item(43, 34).Value = 26  # AH43: 23 -> 26
$ws.Cells.Item(43, 35).Value = 12  # AI43: 13 -> 12
$ws.Cells.Item(43, 36).Value = 7.5  # AJ43: 8 -> 7.5
$ws.Cells.Item(43, 40).Value = 15  # AN43: 17 -> 15
$ws.Cells.Item(43, 41).Value = 26  # AO43: 29 -> 26

# Row 44
$ws.Cells.Item(44, 8).Value = 3.1  # H44: 3.4 -> 3.1
$ws.Cells.Item(44, 9).Value = 2.63  # I44: 2.4 -> 2.63
$ws.Cells.Item(44, 12).Value = 3.2  # L44: 3 -> 3.2
$ws.Cells.Item(44, 19).Value = 1.95  # S44: 1.88 -> 1.95
$ws.Cells.Item(44, 20).Value = 1.9  # T44: 1.98 -> 1.9
$ws.Cells.Item(44, 27).Value = 1.67  # AA44: 1.7 -> 1.67
$ws.Cells.Item(44, 28).Value = 2.1  # AB44: 2.05 -> 2.1
$ws.Cells.Item(44, 35).Value = 10  # AI44: 11 -> 10
$ws.Cells.Item(44, 36).Value = 6  # AJ44: 6.5 -> 6
$ws.Cells.Item(44, 37).Value = 12  # AK44: 13 -> 12
$ws.Cells.Item(44, 39).Value = 151  # AM44: 201 -> 151
$ws.Cells.Item(44, 40).Value = 9.5  # AN44: 9 -> 9.5
$ws.Cells.Item(44, 41).Value = 13  # AO44: 12 -> 13
$ws.Cells.Item(44, 42).Value = 10  # AP44: 9.5 -> 10
$ws.Cells.Item(44, 43).Value = 26  # AQ44: 23 -> 26
$ws.Cells.Item(44, 44).Value = 21  # AR44: 19 -> 21
$ws.Cells.Item(44, 45).Value = 29  # AS44: 26 -> 29

# Row 45
$ws.Cells.Item(45, 7).Value = 1.8  # G45: 1.7 -> 1.8
$ws.Cells.Item(45, 8).Value = 4  # H45: 4.2 -> 4
$ws.Cells.Item(45, 9).Value = 3.9  # I45: 4.33 -> 3.9
$ws.Cells.Item(45, 10).Value = 2.38  # J45: 2.2 -> 2.38
$ws.Cells.Item(45, 11).Value = 2.4  # K45: 2.5 -> 2.4
$ws.Cells.Item(45, 12).Value = 4  # L45: 4.5 -> 4
$ws.Cells.Item(45, 19).Value = 1.57  # S45: 1.5 -> 1.57
$ws.Cells.Item(45, 20).Value = 2.35  # T45: 2.4 -> 2.35
$ws.Cells.Item(45, 21).Value = 1.93  # U45: 1.88 -> 1.93
$ws.Cells.Item(45, 22).Value = 1.93  # V45: 1.98 -> 1.93
$ws.Cells.Item(45, 23).Value = 2.38  # W45: 2.25 -> 2.38
$ws.Cells.Item(45, 24).Value = 1.53  # X45: 1.57 -> 1.53
$ws.Cells.Item(45, 27).Value = 1.53  # AA45: 1.57 -> 1.53
$ws.Cells.Item(45, 28).Value = 2.38  # AB45: 2.25 -> 2.38
$ws.Cells.Item(45, 30).Value = 11  # AD45: 10 -> 11
$ws.Cells.Item(45, 32).Value = 17  # AF45: 15 -> 17
$ws.Cells.Item(45, 33).Value = 13  # AG45: 12 -> 13
$ws.Cells.Item(45, 40).Value = 15  # AN45: 17 -> 15
$ws.Cells.Item(45, 41).Value = 23  # AO45: 26 -> 23
$ws.Cells.Item(45, 42).Value = 13  # AP45: 15 -> 13
$ws.Cells.Item(45, 43).Value = 41  # AQ45: 51 -> 41

# Row 46
$ws.Cells.Item(46, 7).Value = 1.95  # G46: 1.91 -> 1.95
$ws.Cells.Item(46, 8).Value = 3.6  # H46: 3.7 -> 3.6
$ws.Cells.Item(46, 10).Value = 2.6  # J46: 2.5 -> 2.6
$ws.Cells.Item(46, 11).Value = 2.25  # K46: 2.38 -> 2.25
$ws.Cells.Item(46, 13).Value = 1.04  # M46: 1.03 -> 1.04
$ws.Cells.Item(46, 14).Value = 13  # N46: 15 -> 13
$ws.Cells.Item(46, 15).Value = 1.22  # O46: 1.18 -> 1.22
$ws.Cells.Item(46, 16).Value = 4  # P46: 4.5 -> 4
$ws.Cells.Item(46, 19).Value = 1.73  # S46: 1.58 -> 1.73
$ws.Cells.Item(46, 20).Value = 2.08  # T46: 2.25 -> 2.08
$ws.Cells.Item(46, 23).Value = 2.75  # W46: 2.5 -> 2.75
$ws.Cells.Item(46, 24).Value = 1.4  # X46: 1.5 -> 1.4
$ws.Cells.Item(46, 25).Value = 1.33  # Y46: 1.3 -> 1.33
$ws.Cells.Item(46, 26).Value = 3.25  # Z46: 3.4 -> 3.25
$ws.Cells.Item(46, 27).Value = 1.67  # AA46: 1.57 -> 1.67
$ws.Cells.Item(46, 28).Value = 2.1  # AB46: 2.25 -> 2.1
$ws.Cells.Item(46, 29).Value = 9  # AC46: 9.5 -> 9
$ws.Cells.Item(46, 30).Value = 10  # AD46: 11 -> 10
$ws.Cells.Item(46, 33).Value = 15  # AG46: 13 -> 15
$ws.Cells.Item(46, 34).Value = 23  # AH46: 21 -> 23
$ws.Cells.Item(46, 35).Value = 13  # AI46: 15 -> 13
$ws.Cells.Item(46, 37).Value = 13  # AK46: 12 -> 13
$ws.Cells.Item(46, 39).Value = 151  # AM46: 126 -> 151
$ws.Cells.Item(46, 40).Value = 13  # AN46: 15 -> 13
$ws.Cells.Item(46, 44).Value = 29  # AR46: 26 -> 29
$ws.Cells.Item(46, 45).Value = 34  # AS46: 29 -> 34

# Row 47
$ws.Cells.Item(47, 7).Value = 1.75  # G47: 1.73 -> 1.75
$ws.Cells.Item(47, 9).Value = 4.75  # I47: 5 -> 4.75
$ws.Cells.Item(47, 12).Value = 5  # L47: 5.5 -> 5
$ws.Cells.Item(47, 13).Value = 1.05  # M47: 1.06 -> 1.05
$ws.Cells.Item(47, 14).Value = 11  # N47: 9.5 -> 11
$ws.Cells.Item(47, 15).Value = 1.29  # O47: 1.3 -> 1.29
$ws.Cells.Item(47, 16).Value = 3.5  # P47: 3.4 -> 3.5
$ws.Cells.Item(47, 19).Value = 2  # S47: 2.05 -> 2
$ws.Cells.Item(47, 20).Value = 1.85  # T47: 1.8 -> 1.85
$ws.Cells.Item(47, 23).Value = 3.4  # W47: 3.5 -> 3.4
$ws.Cells.Item(47, 24).Value = 1.3  # X47: 1.29 -> 1.3
$ws.Cells.Item(47, 25).Value = 1.4  # Y47: 1.44 -> 1.4
$ws.Cells.Item(47, 26).Value = 2.75  # Z47: 2.63 -> 2.75
$ws.Cells.Item(47, 27).Value = 1.91  # AA47: 1.95 -> 1.91
$ws.Cells.Item(47, 28).Value = 1.91  # AB47: 1.8 -> 1.91
$ws.Cells.Item(47, 29).Value = 7  # AC47: 6.5 -> 7
$ws.Cells.Item(47, 30).Value = 8  # AD47: 7.5 -> 8
$ws.Cells.Item(47, 34).Value = 26  # AH47: 29 -> 26
$ws.Cells.Item(47, 35).Value = 10  # AI47: 9.5 -> 10
$ws.Cells.Item(47, 39).Value = 301  # AM47: 351 -> 301

# Row 48
$ws.Cells.Item(48, 7).Value = 1.91  # G48: 1.85 -> 1.91
$ws.Cells.Item(48, 9).Value = 4  # I48: 4.2 -> 4
$ws.Cells.Item(48, 12).Value = 4  # L48: 4.33 -> 4
$ws.Cells.Item(48, 33).Value = 15  # AG48: 13 -> 15
$ws.Cells.Item(48, 37).Value = 12  # AK48: 13 -> 12
$ws.Cells.Item(48, 40).Value = 13  # AN48: 15 -> 13
$ws.Cells.Item(48, 41).Value = 21  # AO48: 23 -> 21

# Row 56
$ws.Cells.Item(56, 7).Value = 2.62  # G56: 2.45 -> 2.62
$ws.Cells.Item(56, 8).Value = 3.1  # H56: 3.35 -> 3.1
$ws.Cells.Item(56, 9).Value = 2.52  # I56: 2.55 -> 2.52
$ws.Cells.Item(56, 10).Value = 3.25  # J56: 3.05 -> 3.25
$ws.Cells.Item(56, 11).Value = 2.07  # K56: 2.15 -> 2.07
$ws.Cells.Item(56, 13).Value = 1.07  # M56: 1.06 -> 1.07
$ws.Cells.Item(56, 14).Value = 6.7  # N56: 7.2 -> 6.7
$ws.Cells.Item(56, 15).Value = 1.35  # O56: 1.32 -> 1.35
$ws.Cells.Item(56, 16).Value = 2.95  # P56: 3.1 -> 2.95
$ws.Cells.Item(56, 19).Value = 2.02  # S56: 1.93 -> 2.02
$ws.Cells.Item(56, 20).Value = 1.72  # T56: 1.78 -> 1.72
$ws.Cells.Item(56, 23).Value = 3.4  # W56: 3.25 -> 3.4
$ws.Cells.Item(56, 24).Value = 1.27  # X56: 1.3 -> 1.27
$ws.Cells.Item(56, 25).Value = 1.42  # Y56: 1.39 -> 1.42
$ws.Cells.Item(56, 26).Value = 2.67  # Z56: 2.8 -> 2.67
$ws.Cells.Item(56, 29).Value = 8  # AC56: 8.25 -> 8
$ws.Cells.Item(56, 30).Value = 13  # AD56: 12 -> 13
$ws.Cells.Item(56, 31).Value = 10  # AE56: 9.75 -> 10
$ws.Cells.Item(56, 32).Value = 30  # AF56: 25 -> 30
$ws.Cells.Item(56, 33).Value = 23  # AG56: 21 -> 23
$ws.Cells.Item(56, 35).Value = 6.7  # AI56: 7.2 -> 6.7
$ws.Cells.Item(56, 36).Value = 6.1  # AJ56: 6.5 -> 6.1
$ws.Cells.Item(56, 37).Value = 14  # AK56: 14.5 -> 14
$ws.Cells.Item(56, 38).Value = 65  # AL56: 70 -> 65
$ws.Cells.Item(56, 40).Value = 8  # AN56: 8.25 -> 8
$ws.Cells.Item(56, 42).Value = 9.75  # AP56: 10 -> 9.75

# Row 57
$ws.Cells.Item(57, 41).Value = 20  # AO57: 21 -> 20

# Row 58
$ws.Cells.Item(58, 7).Value = 2.95  # G58: 2.67 -> 2.95
$ws.Cells.Item(58, 8).Value = 3.3  # H58: 3.2 -> 3.3
$ws.Cells.Item(58, 9).Value = 2.2  # I58: 2.4 -> 2.2
$ws.Cells.Item(58, 10).Value = 3.6  # J58: 3.3 -> 3.6
$ws.Cells.Item(58, 11).Value = 2.1  # K58: 2.07 -> 2.1
$ws.Cells.Item(58, 12).Value = 2.82  # L58: 3.05 -> 2.82
$ws.Cells.Item(58, 14).Value = 7.4  # N58: 7.3 -> 7.4
$ws.Cells.Item(58, 16).Value = 3.3  # P58: 3.25 -> 3.3
$ws.Cells.Item(58, 19).Value = 1.85  # S58: 1.87 -> 1.85
$ws.Cells.Item(58, 20).Value = 1.85  # T58: 1.83 -> 1.85
$ws.Cells.Item(58, 23).Value = 3  # W58: 3.05 -> 3
$ws.Cells.Item(58, 24).Value = 1.34  # X58: 1.33 -> 1.34
$ws.Cells.Item(58, 25).Value = 1.4  # Y58: 1.42 -> 1.4
$ws.Cells.Item(58, 26).Value = 2.7  # Z58: 2.67 -> 2.7
$ws.Cells.Item(58, 29).Value = 9.5  # AC58: 9.25 -> 9.5
$ws.Cells.Item(58, 30).Value = 15.5  # AD58: 14.5 -> 15.5
$ws.Cells.Item(58, 31).Value = 10.75  # AE58: 10 -> 10.75
$ws.Cells.Item(58, 32).Value = 37  # AF58: 32 -> 37
$ws.Cells.Item(58, 33).Value = 25  # AG58: 22 -> 25
$ws.Cells.Item(58, 34).Value = 32  # AH58: 29 -> 32
$ws.Cells.Item(58, 35).Value = 7.4  # AI58: 7.3 -> 7.4
$ws.Cells.Item(58, 36).Value = 6.5  # AJ58: 6.3 -> 6.5
$ws.Cells.Item(58, 37).Value = 13.5  # AK58: 13 -> 13.5
$ws.Cells.Item(58, 40).Value = 8.25  # AN58: 8.5 -> 8.25
$ws.Cells.Item(58, 41).Value = 11.25  # AO58: 12.5 -> 11.25
$ws.Cells.Item(58, 42).Value = 8.75  # AP58: 9.25 -> 8.75
$ws.Cells.Item(58, 43).Value = 22  # AQ58: 26 -> 22
$ws.Cells.Item(58, 44).Value = 17  # AR58: 19.5 -> 17
$ws.Cells.Item(58, 45).Value = 26  # AS58: 28 -> 26

# Row 59
$ws.Cells.Item(59, 7).Value = 2  # G59: 2.1 -> 2
$ws.Cells.Item(59, 9).Value = 3.35  # I59: 3.1 -> 3.35
$ws.Cells.Item(59, 10).Value = 2.67  # J59: 2.77 -> 2.67
$ws.Cells.Item(59, 12).Value = 4  # L59: 3.8 -> 4
$ws.Cells.Item(59, 16).Value = 2.92  # P59: 2.9 -> 2.92
$ws.Cells.Item(59, 29).Value = 6.8  # AC59: 6.9 -> 6.8
$ws.Cells.Item(59, 30).Value = 9.25  # AD59: 9.5 -> 9.25
$ws.Cells.Item(59, 31).Value = 8.75  # AE59: 9 -> 8.75
$ws.Cells.Item(59, 32).Value = 18  # AF59: 19 -> 18
$ws.Cells.Item(59, 33).Value = 17  # AG59: 18 -> 17
$ws.Cells.Item(59, 40).Value = 9.25  # AN59: 8.75 -> 9.25
$ws.Cells.Item(59, 41).Value = 17  # AO59: 15.5 -> 17
$ws.Cells.Item(59, 42).Value = 12  # AP59: 11.5 -> 12
$ws.Cells.Item(59, 43).Value = 45  # AQ59: 40 -> 45
$ws.Cells.Item(59, 44).Value = 32  # AR59: 30 -> 32
$ws.Cells.Item(59, 45).Value = 45  # AS59: 40 -> 45
